# Fill in four more benchmark names under the existing "9-Grades.xls (modified)"
# entry in column A of Sheet1 (the only sheet with data).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "TRAIL%20INVENTORY%20N#A850A.XLS"
$ws.Range("A4").Value = "2002fairreport.xls"
$ws.Range("A5").Value = "104r.xls"
$ws.Range("A6").Value = "Inventory_Control.xls"

# Widen column A to fit the new (longer) benchmark names, as in the source edit.
$ws.Columns.Item(1).AutoFit() | Out-Null

# The cursor ends up on B6 (one row below the last filled cell) after the edit.
$ws.Range("B6").Select() | Out-Null
